$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number need the Text format
# pre-applied, otherwise Excel (like a real user typing into a General-
# formatted cell) would auto-convert the string into a floating point
# number and lose the exact decimal text (e.g. trailing zeros).
$ws.Range("D2").Value = "67.104.58"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "2.630.07"
$ws.Range("E3").Value = "  -2.04%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.09"
$ws.Range("E5").Value = "  -1.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.68"
$ws.Range("E6").Value = "  +1.41%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -2.38%  "
$ws.Range("D9").Value = "2.629.79"
$ws.Range("E9").Value = "  -2.07%  "
$ws.Range("E10").Value = "  -1.47%  "
$ws.Range("E11").Value = "  +1.54%  "
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("E13").Value = "  +0.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.74"
$ws.Range("E14").Value = "  -0.48%  "
$ws.Range("E15").Value = "  -2.07%  "
$ws.Range("E16").Value = "  -1.52%  "
$ws.Range("D17").Value = "66.879.09"
$ws.Range("E17").Value = "  -1.15%  "
$ws.Range("D18").Value = "2.629.32"
$ws.Range("E18").Value = "  -2.27%  "
$ws.Range("E19").Value = "  +4.76%  "
$ws.Range("E20").Value = "  +7.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "357.36"
$ws.Range("E21").Value = "  -2.14%  "
$ws.Range("E22").Value = "  -1.97%  "
$ws.Range("E23").Value = "  -3.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.92"
$ws.Range("E24").Value = "  +10.44%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("E26").Value = "  -4.81%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "70.22"
$ws.Range("E27").Value = "  -3.38%  "
$ws.Range("D28").Value = "2.757.65"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("E30").Value = "  -1.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "550.25"
$ws.Range("E31").Value = "  -2.29%  "
$ws.Range("E32").Value = "  -0.59%  "
$ws.Range("E33").Value = "  -1.77%  "
$ws.Range("E34").Value = "  -2.02%  "
$ws.Range("E35").Value = "  +5.21%  "
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("E37").Value = "  -4.72%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "155.48"
$ws.Range("E38").Value = "  -0.61%  "
$ws.Range("E39").Value = "  -2.63%  "
$ws.Range("E40").Value = "  -1.67%  "
$ws.Range("E41").Value = "  -2.07%  "
$ws.Range("E42").Value = "  -2.22%  "
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.43"
$ws.Range("E45").Value = "  -3.86%  "
$ws.Range("E46").Value = "  -0.66%  "
$ws.Range("E47").Value = "  -1.71%  "
$ws.Range("E48").Value = "  -1.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "151.60"
$ws.Range("E49").Value = "  -1.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.78"
$ws.Range("E50").Value = "  -1.64%  "
$ws.Range("E51").Value = "  -0.77%  "
